$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that are being dropped from the output (by their
# original row numbers, bottom-most first so earlier row numbers stay valid).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(3).Delete()

# Remove the old "max" column (old column C); this shifts the old D/E
# "prediction"/"rejection-f" columns left into C/D, matching the target.
$ws.Columns.Item(3).Delete()

# Update the "B" (1-s__CAG-345 sp000433315) numeric column with the new
# per-MAG values.
$ws.Range("B2").Value = 1.693434838499416
$ws.Range("B3").Value = 1.463256242741473
$ws.Range("B4").Value = 1.748763174246024
$ws.Range("B5").Value = 1.693434838499416
$ws.Range("B6").Value = 1.68482165972379
$ws.Range("B7").Value = 1.314242988281904
